$d = $word.ActiveDocument

# --------------------------------------------------------------------
# Original paragraph layout (1-based, verified against before.docx):
#   1  Questions
#   2  Does the religious affiliation of an MP ...
#   3  If the dominant religion ... Christianity ... being a Christian?
#   4  If the dominant religion ... Islam ... being a Muslim?
#   5  If the dominant religion ... Traditional Religion ... Traditionalist?
#   6  What is the dominant religion of MPs in Ghana?
#   7  As a campaign Team Advisor of any of the political parties ...
#   8  (empty paragraph)
#   9  Results
#  10  89.82 percent of the MPs have the same dominant religion ...
#  11  If the dominant religion ... Christianity ... being a Christian?
#  12  If the dominant religion ... Islam ... Muslim is ……………………..
#  13  If the dominant religion ... Traditional Religion ... is ………………………………
#  14  T|he do|minant religion of MPs in Ghana is …………………………..
#  15  From the results, as a campaign Team advisor ...
#
# Target keeps only: 1, 2, 6, 8, 9, 10, 15, and a rewritten 14 -- i.e.
# paragraphs 3, 4, 5, 7, 11, 12, 13 are dropped outright.
# --------------------------------------------------------------------

function Assert-ParagraphStartsWith($doc, $index, $expectedPrefix) {
    $actual = $doc.Paragraphs($index).Range.Text
    if ($actual.Length -lt $expectedPrefix.Length -or
        $actual.Substring(0, $expectedPrefix.Length) -ne $expectedPrefix) {
        throw "Paragraph $index does not start with expected text. Expected prefix: '$expectedPrefix' Actual: '$actual'"
    }
}

# Sanity-check the paragraphs we are about to delete, then delete them
# bottom-up so the indices above remain valid through the whole pass.
Assert-ParagraphStartsWith $d 13 "If the dominant religion in a particular constituency is Traditional Religion, the likelihood"
$d.Paragraphs(13).Range.Delete()

Assert-ParagraphStartsWith $d 12 "If the dominant religion in a particular constituency is Islam, the likelihood"
$d.Paragraphs(12).Range.Delete()

Assert-ParagraphStartsWith $d 11 "If the dominant religion in a particular constituency is Christianity, what is the likelihood"
$d.Paragraphs(11).Range.Delete()

Assert-ParagraphStartsWith $d 7 "As a campaign Team Advisor of any of the political parties"
$d.Paragraphs(7).Range.Delete()

Assert-ParagraphStartsWith $d 5 "If the dominant religion in a particular constituency is "
$d.Paragraphs(5).Range.Delete()

Assert-ParagraphStartsWith $d 4 "If the dominant religion "
$d.Paragraphs(4).Range.Delete()

Assert-ParagraphStartsWith $d 3 "If the dominant religion "
$d.Paragraphs(3).Range.Delete()

# --------------------------------------------------------------------
# Rewrite the mangled "T|he do|minant religion of MPs in Ghana is
# ………………………….." paragraph (now index 7) as a clean sentence ending
# in "Christianity".
# --------------------------------------------------------------------
Assert-ParagraphStartsWith $d 7 "T"

$p = $d.Paragraphs(7)
$pStart = $p.Range.Start
$pEnd = $p.Range.End   # Range includes the trailing paragraph mark at End-1

# Drop everything after the leading "T" (but keep the paragraph mark).
$tail = $d.Range($pStart + 1, $pEnd - 1)
$tail.Delete()

# Turn the lone "T" run into the full lead-in sentence.
$lead = $d.Range($pStart, $pStart + 1)
$lead.Text = "The dominant religion of MPs in Ghana is "

# Append "Christianity" right before the paragraph mark.
$p = $d.Paragraphs(7)
$insertAt = $d.Range($p.Range.End - 1, $p.Range.End - 1)
$insertAt.InsertAfter("Christianity")

$expectedFinal = "The dominant religion of MPs in Ghana is Christianity"
$final = $d.Paragraphs(7).Range.Text
if ($final.Substring(0, $final.Length - 1) -ne $expectedFinal) {
    throw "Unexpected final paragraph text: '$final'"
}
